$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric (single decimal point) would be
# auto-converted from text to a Number by Excels input parser, which
# both changes the cell type and can introduce floating-point drift
# (e.g. "305.29" -> 305.29000000000002). Force those specific cells to
# the Text number format first so the literal string is preserved.
$ws.Range("D2").Value = '42.230.18'
$ws.Range("E2").Value = '  +0.84%  '
$ws.Range("D3").Value = '2.264.73'
$ws.Range("E3").Value = '  -0.30%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.29'
$ws.Range("E5").Value = '  +0.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '96.93'
$ws.Range("E6").Value = '  +4.28%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.529'
$ws.Range("E7").Value = '  -0.39%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.489'
$ws.Range("E9").Value = '  +0.63%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.53'
$ws.Range("E10").Value = '  +8.41%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0794'
$ws.Range("E11").Value = '  -0.35%  '
$ws.Range("E12").Value = '  -1.11%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.64'
$ws.Range("E13").Value = '  -1.00%  '
$ws.Range("D14").Value = '2.624.26'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.32'
$ws.Range("E15").Value = '  +0.02%  '
$ws.Range("D16").Value = '2.267.47'
$ws.Range("E16").Value = '  -1.82%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.792'
$ws.Range("E17").Value = '  +1.70%  '
$ws.Range("D18").Value = '42.143.17'
$ws.Range("E18").Value = '  +0.84%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.47'
$ws.Range("E19").Value = '  -2.21%  '
$ws.Range("D20").Value = '0.0₃0908'
$ws.Range("E20").Value = '  -0.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.95'
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.42'
$ws.Range("E22").Value = '  +0.35%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.66'
$ws.Range("E23").Value = '  -2.79%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.58'
$ws.Range("E24").Value = '  -0.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.95'
$ws.Range("E25").Value = '  +0.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  -0.38%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.77'
$ws.Range("E27").Value = '  -1.27%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.39'
$ws.Range("E28").Value = '  +5.83%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.49'
$ws.Range("E29").Value = '  -0.48%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.11'
$ws.Range("E30").Value = '  +1.87%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '160.27'
$ws.Range("E31").Value = '  -0.48%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.24'
$ws.Range("E32").Value = '  -0.40%  '
$ws.Range("E33").Value = '  +0.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.14'
$ws.Range("E34").Value = '  +3.94%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0739'
$ws.Range("E35").Value = '  -0.74%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.04'
$ws.Range("E36").Value = '  -0.58%  '
$ws.Range("E37").Value = '  -1.36%  '
$ws.Range("E38").Value = '  -0.33%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.82'
$ws.Range("E39").Value = '  +0.97%  '
$ws.Range("E40").Value = '  -1.66%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.05'
$ws.Range("E41").Value = '  +2.89%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.42'
$ws.Range("E42").Value = '  +13.98%  '
$ws.Range("D43").Value = '1.991.84'
$ws.Range("E43").Value = '  -0.75%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0288'
$ws.Range("E44").Value = '  +1.96%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.70'
$ws.Range("E45").Value = '  -4.48%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.95'
$ws.Range("E46").Value = '  -3.60%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.93'
$ws.Range("E47").Value = '  +0.47%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '53.19'
$ws.Range("E48").Value = '  +0.42%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.51'
$ws.Range("E49").Value = '  -0.49%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '71.88'
$ws.Range("E50").Value = '  -0.59%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '91.19'
$ws.Range("E51").Value = '  -0.26%  '
